# "Docu actualizada con las pruebas"
# The student filled in the test-case results table (sheet "Pruebas"):
#   - Marked test cases 1 through 29 (rows 64-92) as "OK" in the Estado column (E).
#   - Fixed a typo in the description of [Prueba29] (row 92, column I) where the
#     opening bracket "[" was missing.
#   - Updated the view (zoom level and selected cell) to reflect where they left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pruebas")

# Mark test cases 1-29 (rows 64 to 92) as "OK"
$ws.Range("E64:E92").Value = "OK"

# Fix the typo in the description of Prueba29: add the missing leading "["
$ws.Range("I92").Value = "[Prueba29] Identificarse en la aplicación y enviar un mensaje a un amigo, validar que el mensaje enviado aparece en el chat."

# Update the sheet view: zoom out to 55% and select E92
$ws.Activate()
$excel.ActiveWindow.Zoom = 55
$ws.Range("E92").Select()
